$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4 (new): previous row-3 data ("$/caja 12 kilos granel", O'Higgins) moves here.
$ws.Cells.Item(4, 1).Value = 11
$ws.Cells.Item(4, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(4, 3).Value = "Bíobío"
$ws.Cells.Item(4, 4).Value = 44334
$ws.Cells.Item(4, 5).Value = 8
$ws.Cells.Item(4, 6).Value = "Fruta"
$ws.Cells.Item(4, 7).Value = 100107
$ws.Cells.Item(4, 8).Value = "Otros"
$ws.Cells.Item(4, 9).Value = 100107001
$ws.Cells.Item(4, 10).Value = "Caqui"
$ws.Cells.Item(4, 11).Value = "Mankaki"
$ws.Cells.Item(4, 12).Value = "Primera"
$ws.Cells.Item(4, 13).Value = 100
$ws.Cells.Item(4, 14).Value = 11000
$ws.Cells.Item(4, 15).Value = 12000
$ws.Cells.Item(4, 16).Value = 11500
$ws.Cells.Item(4, 17).Value = "`$/caja 12 kilos granel"
$ws.Cells.Item(4, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(4, 19).Value = 11500
$ws.Cells.Item(4, 20).Value = 1
$ws.Cells.Item(4, 4).NumberFormat = $ws.Cells.Item(2, 4).NumberFormat

# Row 3: previous row-2 data ("$/caja 18 kilos granel", Curicó) moves here.
$ws.Cells.Item(3, 1).Value = 11
$ws.Cells.Item(3, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(3, 3).Value = "Bíobío"
$ws.Cells.Item(3, 4).Value = 44330
$ws.Cells.Item(3, 5).Value = 8
$ws.Cells.Item(3, 6).Value = "Fruta"
$ws.Cells.Item(3, 7).Value = 100107
$ws.Cells.Item(3, 8).Value = "Otros"
$ws.Cells.Item(3, 9).Value = 100107001
$ws.Cells.Item(3, 10).Value = "Caqui"
$ws.Cells.Item(3, 11).Value = "Mankaki"
$ws.Cells.Item(3, 12).Value = "Primera"
$ws.Cells.Item(3, 13).Value = 100
$ws.Cells.Item(3, 14).Value = 15000
$ws.Cells.Item(3, 15).Value = 16000
$ws.Cells.Item(3, 16).Value = 15500
$ws.Cells.Item(3, 17).Value = "`$/caja 18 kilos granel"
$ws.Cells.Item(3, 18).Value = "Provincia de Curicó"
$ws.Cells.Item(3, 19).Value = 861
$ws.Cells.Item(3, 20).Value = 18
$ws.Cells.Item(3, 4).NumberFormat = $ws.Cells.Item(2, 4).NumberFormat

# Row 2 (new content): "$/caja 12 kilos empedrada" entry.
$ws.Cells.Item(2, 1).Value = 11
$ws.Cells.Item(2, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(2, 3).Value = "Bíobío"
$ws.Cells.Item(2, 4).Value = 44707
$ws.Cells.Item(2, 5).Value = 8
$ws.Cells.Item(2, 6).Value = "Fruta"
$ws.Cells.Item(2, 7).Value = 100107
$ws.Cells.Item(2, 8).Value = "Otros"
$ws.Cells.Item(2, 9).Value = 100107001
$ws.Cells.Item(2, 10).Value = "Caqui"
$ws.Cells.Item(2, 11).Value = "Mankaki"
$ws.Cells.Item(2, 12).Value = "Primera"
$ws.Cells.Item(2, 13).Value = 60
$ws.Cells.Item(2, 14).Value = 12000
$ws.Cells.Item(2, 15).Value = 13000
$ws.Cells.Item(2, 16).Value = 12500
$ws.Cells.Item(2, 17).Value = "`$/caja 12 kilos empedrada"
$ws.Cells.Item(2, 18).Value = "Provincia de Curicó"
$ws.Cells.Item(2, 19).Value = 1042
$ws.Cells.Item(2, 20).Value = 12
